$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 2611
$ws.Cells.Item(4, 6).Value = 22
$ws.Cells.Item(5, 6).Value = 1314
$ws.Cells.Item(7, 6).Value = 3195
$ws.Cells.Item(8, 6).Value = 381
$ws.Cells.Item(9, 6).Value = 175
$ws.Cells.Item(10, 6).Value = 48
$ws.Cells.Item(11, 6).Value = 8124
$ws.Cells.Item(12, 6).Value = 422
$ws.Cells.Item(13, 6).Value = 76
$ws.Cells.Item(15, 6).Value = 22
$ws.Cells.Item(16, 6).Value = 291
$ws.Cells.Item(18, 6).Value = 59
$ws.Cells.Item(20, 6).Value = 311
$ws.Cells.Item(21, 6).Value = 10137
$ws.Cells.Item(22, 6).Value = 30
$ws.Cells.Item(24, 6).Value = 35
$ws.Cells.Item(26, 6).Value = 381
$ws.Cells.Item(30, 6).Value = 82
$ws.Cells.Item(32, 6).Value = 43
$ws.Cells.Item(33, 6).Value = 2060
$ws.Cells.Item(34, 6).Value = 30
$ws.Cells.Item(35, 6).Value = 30
$ws.Cells.Item(36, 6).Value = 2100
$ws.Cells.Item(37, 6).Value = 4030
$ws.Cells.Item(38, 6).Value = 254
$ws.Cells.Item(40, 6).Value = 2174
$ws.Cells.Item(41, 6).Value = 1215
$ws.Cells.Item(42, 6).Value = 144
$ws.Cells.Item(43, 6).Value = 297
$ws.Cells.Item(44, 6).Value = 222
$ws.Cells.Item(45, 6).Value = 28
$ws.Cells.Item(46, 6).Value = 89
$ws.Cells.Item(47, 6).Value = 80
$ws.Cells.Item(49, 6).Value = 52

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 9
$ws.Cells.Item(5, 6).Value = 170
$ws.Cells.Item(6, 6).Value = 38
$ws.Cells.Item(7, 6).Value = 6
$ws.Cells.Item(13, 6).Value = 27
$ws.Cells.Item(16, 6).Value = 172

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 3

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 2611
$ws.Cells.Item(4, 6).Value = 170
$ws.Cells.Item(5, 6).Value = 22
$ws.Cells.Item(6, 6).Value = 1314
$ws.Cells.Item(8, 6).Value = 3195
$ws.Cells.Item(9, 6).Value = 381
$ws.Cells.Item(11, 6).Value = 175
$ws.Cells.Item(12, 6).Value = 48
$ws.Cells.Item(13, 6).Value = 8124
$ws.Cells.Item(14, 6).Value = 422
$ws.Cells.Item(15, 6).Value = 76
$ws.Cells.Item(17, 6).Value = 22
$ws.Cells.Item(18, 6).Value = 291
$ws.Cells.Item(19, 6).Value = 59
$ws.Cells.Item(21, 6).Value = 311
$ws.Cells.Item(22, 6).Value = 10137
$ws.Cells.Item(24, 6).Value = 35
$ws.Cells.Item(26, 6).Value = 381
$ws.Cells.Item(28, 6).Value = 27
$ws.Cells.Item(30, 6).Value = 82
$ws.Cells.Item(32, 6).Value = 43
$ws.Cells.Item(33, 6).Value = 2060
$ws.Cells.Item(34, 6).Value = 30
$ws.Cells.Item(35, 6).Value = 2100
$ws.Cells.Item(36, 6).Value = 4030
$ws.Cells.Item(37, 6).Value = 254
$ws.Cells.Item(39, 6).Value = 2175
$ws.Cells.Item(41, 6).Value = 1215
$ws.Cells.Item(42, 6).Value = 144
$ws.Cells.Item(43, 6).Value = 297
$ws.Cells.Item(44, 6).Value = 222
$ws.Cells.Item(45, 6).Value = 28
$ws.Cells.Item(46, 6).Value = 89
$ws.Cells.Item(47, 6).Value = 80
$ws.Cells.Item(49, 6).Value = 52
